# Update "想去人数" (want-to-go count) values in column F
# for worksheets "展览" (sheet1) and "全部类型" (sheet4),
# matching the output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new F value for sheet "展览"
$exhibitUpdates = @{
    3  = 481
    7  = 7358
    10 = 3072
    16 = 24
    17 = 745
    20 = 196
    23 = 121
    24 = 361
    27 = 72
    28 = 114
    29 = 2099
    30 = 614
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new F value for sheet "全部类型"
$allUpdates = @{
    4  = 481
    8  = 7358
    12 = 3072
    19 = 24
    23 = 745
    26 = 196
    32 = 121
    33 = 361
    36 = 72
    37 = 114
    38 = 2099
    39 = 614
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
